$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Phone Number"), shifting existing Purchase_Type/Payment_Mode
# columns (D,E) over to (E,F).
$ws.Columns("D").Insert()

# Header row
$ws.Range("D1").Value = "Phone Number"

# Data rows - prefix with an apostrophe so the numeric-looking phone numbers are
# stored as text (shared strings) rather than being coerced to numbers.
$ws.Range("D2").Value = "'7016763640"
$ws.Range("D3").Value = "'7405802474"
